# Correlaciones en SubSaturado y Saturado
# Corrige los valores de Bg (columna E) de la hoja, que estaban expresados
# con una escala 1000 veces mayor de lo correcto.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.00432355
$ws.Range("E3").Value = 0.00443585
$ws.Range("E4").Value = 0.00488505
$ws.Range("E5").Value = 0.00555885

# Restablece la celda activa/seleccion tal como quedo tras la edicion
$ws.Range("H6").Select()
